$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 122, shifting existing rows
# 122-186 down to 124-188 (preserves all their original values/formats).
$ws.Range("A122:A123").EntireRow.Insert()

# Populate new row 122 (new weekly data point, date 2022-01-11 / serial 44572)
$ws.Cells.Item(122, 1).Value = 11
$ws.Cells.Item(122, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(122, 3).Value = "Bíobío"
$ws.Cells.Item(122, 4).Value = 44572
$ws.Cells.Item(122, 5).Value = 8
$ws.Cells.Item(122, 6).Value = 100112017
$ws.Cells.Item(122, 7).Value = "Apio"
$ws.Cells.Item(122, 8).Value = "Americana (o)"
$ws.Cells.Item(122, 9).Value = "Primera"
$ws.Cells.Item(122, 10).Value = 350
$ws.Cells.Item(122, 11).Value = 6500
$ws.Cells.Item(122, 12).Value = 7000
$ws.Cells.Item(122, 13).Value = 6714
$ws.Cells.Item(122, 14).Value = "$/docena de matas"
$ws.Cells.Item(122, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(122, 16).Value = 1119
$ws.Cells.Item(122, 17).Value = 6
$ws.Cells.Item(122, 18).Value = "Hortaliza"

# Populate new row 123 (same date, "Segunda" quality)
$ws.Cells.Item(123, 1).Value = 11
$ws.Cells.Item(123, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(123, 3).Value = "Bíobío"
$ws.Cells.Item(123, 4).Value = 44572
$ws.Cells.Item(123, 5).Value = 8
$ws.Cells.Item(123, 6).Value = 100112017
$ws.Cells.Item(123, 7).Value = "Apio"
$ws.Cells.Item(123, 8).Value = "Americana (o)"
$ws.Cells.Item(123, 9).Value = "Segunda"
$ws.Cells.Item(123, 10).Value = 260
$ws.Cells.Item(123, 11).Value = 5000
$ws.Cells.Item(123, 12).Value = 5500
$ws.Cells.Item(123, 13).Value = 5269
$ws.Cells.Item(123, 14).Value = "$/docena de matas"
$ws.Cells.Item(123, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(123, 16).Value = 878
$ws.Cells.Item(123, 17).Value = 6
$ws.Cells.Item(123, 18).Value = "Hortaliza"
